{"js": "// Remove the \"Consider how Buede...\" prompt paragraph and the two\n// Heading1 paragraphs (\"Grady's Position\" / \"Buede's Position\") that\n// followed it, now that the assignment has been completed. The\n// paragraph immediately preceding them (\"Steve Mazza\") now runs\n// directly into the \"Resolution\" heading.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet startIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim().indexOf(\"Consider how\") === 0) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // Delete the prompt paragraph plus the two heading paragraphs right\n  // after it (\"Grady's Position\" and \"Buede's Position\").\n  const countToRemove = 3;\n  for (let i = startIndex; i < startIndex + countToRemove && i < paragraphs.items.length; i++) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"Consider how Buede...\" prompt paragraph and the two\n# Heading1 paragraphs (\"Grady's Position\" / \"Buede's Position\") that\n# followed it, now that the assignment has been completed. The\n# paragraph immediately preceding them (\"Steve Mazza\") now runs\n# directly into the \"Resolution\" heading.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.Trim()\n    if ($text.StartsWith(\"Consider how\")) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $startPara = $target\n    $endPara = $startPara.Next().Next()\n    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $range.Delete()\n}\n"}
